$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 12666.667
$ws.Range("I21").Value = 0
$ws.Range("J21").Value = 12666.667
$ws.Range("K21").Value = 0
$ws.Range("L21").Value = 12666.667
$ws.Range("N21").Value = -13602.667
$ws.Range("H23").Value = 12666.667
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 12666.667
$ws.Range("K23").Value = 0
$ws.Range("L23").Value = 12666.667
$ws.Range("N23").Value = -13134.667
$ws.Range("H106").Value = 1936
$ws.Range("I106").Value = 1992.6154
$ws.Range("J106").Value = 1200
$ws.Range("K106").Value = 1992.6154
$ws.Range("L106").Value = 1200
$ws.Range("M106").Value = -1361.6154
$ws.Range("N106").Value = -2462
$ws.Range("M21").ClearContents()
$ws.Range("M23").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1785.1915
$ws.Range("I61").Value = 1032.6111
$ws.Range("J61").Value = 2252.3103
$ws.Range("K61").Value = 1032.6111
$ws.Range("L61").Value = 2252.3103
$ws.Range("M61").Value = -820.6111000000001
$ws.Range("N61").Value = -2676.3103
$ws.Range("H74").Value = 1588.3062
$ws.Range("I74").Value = 928.52
$ws.Range("J74").Value = 2275.5833
$ws.Range("K74").Value = 928.52
$ws.Range("L74").Value = 2275.5833
$ws.Range("M74").Value = -54.51999999999998
$ws.Range("N74").Value = -4023.5833
$ws.Range("H77").Value = 1588.3062
$ws.Range("I77").Value = 928.52
$ws.Range("J77").Value = 2275.5833
$ws.Range("K77").Value = 4642.6
$ws.Range("L77").Value = 11377.9165
$ws.Range("M77").Value = -274.6000000000004
$ws.Range("N77").Value = -20113.9165
$ws.Range("H136").Value = 1785.1915
$ws.Range("I136").Value = 1032.6111
$ws.Range("J136").Value = 2252.3103
$ws.Range("K136").Value = 3097.8333
$ws.Range("L136").Value = 6756.9309
$ws.Range("M136").Value = -547.8333000000002
$ws.Range("N136").Value = -11856.9309

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H108").Value = 36000
$ws.Range("J108").Value = 36000
$ws.Range("L108").Value = 36000
$ws.Range("N108").Value = -43680

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 15326.054
$ws.Range("I31").Value = 39562.04
$ws.Range("K31").Value = 39562.04
$ws.Range("M31").Value = -39267.04
$ws.Range("H34").Value = 15326.054
$ws.Range("I34").Value = 39562.04
$ws.Range("K34").Value = 39562.04
$ws.Range("M34").Value = -39360.04
$ws.Range("H36").Value = 7933.3335
$ws.Range("I36").Value = 9800
$ws.Range("J36").Value = 7000
$ws.Range("K36").Value = 9800
$ws.Range("L36").Value = 7000
$ws.Range("M36").Value = -9412
$ws.Range("N36").Value = -7776
$ws.Range("H40").Value = 7933.3335
$ws.Range("I40").Value = 9800
$ws.Range("J40").Value = 7000
$ws.Range("K40").Value = 9800
$ws.Range("L40").Value = 7000
$ws.Range("M40").Value = -9640
$ws.Range("N40").Value = -7320
$ws.Range("H58").Value = 17382.5
$ws.Range("I58").Value = 1716.9231
$ws.Range("J58").Value = 85266.664
$ws.Range("K58").Value = 1716.9231
$ws.Range("L58").Value = 85266.664
$ws.Range("M58").Value = -1513.9231
$ws.Range("N58").Value = -85672.664
$ws.Range("H96").Value = 17225
$ws.Range("J96").Value = 17225
$ws.Range("L96").Value = 17225
$ws.Range("N96").Value = -22717
$ws.Range("H132").Value = 2809.2222
$ws.Range("I132").Value = 2676.238
$ws.Range("J132").Value = 3274.6667
$ws.Range("K132").Value = 8028.714
$ws.Range("L132").Value = 9824.000100000001
$ws.Range("M132").Value = -5498.714
$ws.Range("N132").Value = -14884.0001
$ws.Range("H136").Value = 17382.5
$ws.Range("I136").Value = 1716.9231
$ws.Range("J136").Value = 85266.664
$ws.Range("K136").Value = 5150.7693
$ws.Range("L136").Value = 255799.992
$ws.Range("M136").Value = -2600.7693
$ws.Range("N136").Value = -260899.992

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 5643.175
$ws.Range("I5").Value = 512.0571
$ws.Range("J5").Value = 41561
$ws.Range("K5").Value = 1536.1713
$ws.Range("L5").Value = 124683
$ws.Range("M5").Value = -1424.1713
$ws.Range("N5").Value = -124907
$ws.Range("H75").Value = 10568.583
$ws.Range("I75").Value = 941.6667
$ws.Range("J75").Value = 13777.556
$ws.Range("K75").Value = 2825.0001
$ws.Range("L75").Value = 41332.66800000001
$ws.Range("M75").Value = -1827.0001
$ws.Range("N75").Value = -43328.66800000001
$ws.Range("H78").Value = 10568.583
$ws.Range("I78").Value = 941.6667
$ws.Range("J78").Value = 13777.556
$ws.Range("K78").Value = 8475.0003
$ws.Range("L78").Value = 123998.004
$ws.Range("M78").Value = -3483.0003
$ws.Range("N78").Value = -133982.004
$ws.Range("H135").Value = 5643.175
$ws.Range("I135").Value = 512.0571
$ws.Range("J135").Value = 41561
$ws.Range("K135").Value = 4608.5139
$ws.Range("L135").Value = 374049
$ws.Range("M135").Value = -2073.5139
$ws.Range("N135").Value = -379119
$ws.Range("H137").Value = 2935.6667
$ws.Range("I137").Value = 2363
$ws.Range("J137").Value = 3590.1428
$ws.Range("K137").Value = 7089
$ws.Range("L137").Value = 10770.4284
$ws.Range("M137").Value = -1989
$ws.Range("N137").Value = -20970.4284

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H52").Value = 11648.167
$ws.Range("I52").Value = 6030
$ws.Range("J52").Value = 12771.8
$ws.Range("K52").Value = 6030
$ws.Range("L52").Value = 12771.8
$ws.Range("M52").Value = -5771
$ws.Range("N52").Value = -13289.8
$ws.Range("H97").Value = 125002980
$ws.Range("I97").Value = 200003400
$ws.Range("J97").Value = 2266.3333
$ws.Range("K97").Value = 200003400
$ws.Range("L97").Value = 2266.3333
$ws.Range("M97").Value = -200002904
$ws.Range("N97").Value = -3258.3333
$ws.Range("H136").Value = 47342.445
$ws.Range("J136").Value = 47342.445
$ws.Range("L136").Value = 142027.335
$ws.Range("N136").Value = -147127.335

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 0
$ws.Range("I54").Value = 0
$ws.Range("J54").Value = 0
$ws.Range("K54").Value = 0
$ws.Range("L54").Value = 0
$ws.Range("H97").Value = 30571
$ws.Range("J97").Value = 30571
$ws.Range("L97").Value = 30571
$ws.Range("N97").Value = -32553
$ws.Range("M54").ClearContents()
$ws.Range("N54").ClearContents()
